$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(135).Insert()

$ws.Cells.Item(135, 1).Value = 11
$ws.Cells.Item(135, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(135, 3).Value = "Bíobío"
$ws.Cells.Item(135, 4).Value = 45006
$ws.Cells.Item(135, 5).Value = 8
$ws.Cells.Item(135, 6).Value = 100112043
$ws.Cells.Item(135, 7).Value = "Pepino ensalada"
$ws.Cells.Item(135, 8).Value = "Sin especificar"
$ws.Cells.Item(135, 9).Value = "Primera"
$ws.Cells.Item(135, 10).Value = 100
$ws.Cells.Item(135, 11).Value = 7000
$ws.Cells.Item(135, 12).Value = 7500
$ws.Cells.Item(135, 13).Value = 7250
$ws.Cells.Item(135, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(135, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(135, 16).Value = 121
$ws.Cells.Item(135, 17).Value = 60
$ws.Cells.Item(135, 18).Value = "Hortaliza"
